$wb = $excel.ActiveWorkbook

# --- m_steam sheet: header renamed from "m steam [kg/h]" to "m steam [kg/s]"
#     and data converted from kg/h to kg/s (divide by 3600).
$wsSteam = $wb.Worksheets.Item("m_steam")
$wsSteam.Range("A1").Value = "m steam [kg/s]"
for ($r = 4; $r -le 10; $r++) {
    for ($c = 3; $c -le 8; $c++) {
        $cell = $wsSteam.Cells.Item($r, $c)
        $cell.Value = $cell.Value2 / 3600
    }
}

# --- m_fuel sheet: header renamed from "m fuel [kg/h]" to "m fuel [kg/s]"
#     and data converted from kg/h to kg/s (divide by 3600).
$wsFuel = $wb.Worksheets.Item("m_fuel")
$wsFuel.Range("A1").Value = "m fuel [kg/s]"
for ($r = 4; $r -le 10; $r++) {
    for ($c = 3; $c -le 8; $c++) {
        $cell = $wsFuel.Cells.Item($r, $c)
        $cell.Value = $cell.Value2 / 3600
    }
}

# --- Restore the active selections on each sheet view.
$wsSteam.Activate()
$wsSteam.Range("O15").Select()

$wsFuel.Activate()
$wsFuel.Range("J4:O10").Select()
